# chore: update Sheets via scheduled runner
#
# Refreshes cached market-board derived figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ -> columns H..N) for a handful of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 3426.037
$ws.Range("I74").Value = 3161.2778
$ws.Range("K74").Value = 3161.2778
$ws.Range("M74").Value = -2225.2778

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 3426.037
$ws.Range("I77").Value = 3161.2778
$ws.Range("K77").Value = 15806.389
$ws.Range("M77").Value = -11126.389

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1503.05
$ws.Range("I137").Value = 1032.9286
$ws.Range("J137").Value = 2600
$ws.Range("K137").Value = 3098.7858
$ws.Range("L137").Value = 7800
$ws.Range("M137").Value = -548.7857999999997
$ws.Range("N137").Value = -12900

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1808.8918
$ws.Range("I138").Value = 1582.4584
$ws.Range("J138").Value = 2226.923
$ws.Range("K138").Value = 4747.3752
$ws.Range("L138").Value = 6680.768999999999
$ws.Range("M138").Value = 392.6247999999996
$ws.Range("N138").Value = -16960.769

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 20833996
$ws.Range("I5").Value = 41667104
$ws.Range("J5").Value = 887.5
$ws.Range("K5").Value = 41667104
$ws.Range("L5").Value = 887.5
$ws.Range("M5").Value = -41666992
$ws.Range("N5").Value = -1111.5

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 777.2069
$ws.Range("I74").Value = 689.125
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 689.125
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 184.875
$ws.Range("N74").Value = -2948

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 777.2069
$ws.Range("I77").Value = 689.125
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 3445.625
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = 922.375
$ws.Range("N77").Value = -14736

# Row 102: Smells of Rich Tama-hagane (N102 previously absent -> now populated)
$ws.Range("H102").Value = 1675
$ws.Range("I102").Value = 1400
$ws.Range("J102").Value = 1950
$ws.Range("K102").Value = 1400
$ws.Range("L102").Value = 1950
$ws.Range("M102").Value = 222
$ws.Range("N102").Value = -5194

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 20833996
$ws.Range("I4").Value = 41667104
$ws.Range("J4").Value = 887.5
$ws.Range("K4").Value = 41667104
$ws.Range("L4").Value = 887.5
$ws.Range("M4").Value = -41666989
$ws.Range("N4").Value = -1117.5

# Row 94: High Steal
$ws.Range("H94").Value = 446.42856
$ws.Range("I94").Value = 496.64706
$ws.Range("J94").Value = 368.81818
$ws.Range("K94").Value = 496.64706
$ws.Range("L94").Value = 368.81818
$ws.Range("M94").Value = -45.64706000000001
$ws.Range("N94").Value = -1270.81818

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move (L15 removed, M15 now carries the old L15-derived value)
$ws.Range("H15").Value = 70009
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 70009
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 70009
$ws.Range("N15").Value = -70349

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2605998.5
$ws.Range("I31").Value = 1290.4348
$ws.Range("K31").Value = 1290.4348
$ws.Range("M31").Value = -995.4348

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2605998.5
$ws.Range("I34").Value = 1290.4348
$ws.Range("K34").Value = 1290.4348
$ws.Range("M34").Value = -1088.4348

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1807.826
$ws.Range("I58").Value = 1223.75
$ws.Range("J58").Value = 3142.8572
$ws.Range("K58").Value = 1223.75
$ws.Range("L58").Value = 3142.8572
$ws.Range("M58").Value = -1020.75
$ws.Range("N58").Value = -3548.8572

# Row 136: Turali Quality
$ws.Range("H136").Value = 1807.826
$ws.Range("I136").Value = 1223.75
$ws.Range("J136").Value = 3142.8572
$ws.Range("K136").Value = 3671.25
$ws.Range("L136").Value = 9428.571599999999
$ws.Range("M136").Value = -1121.25
$ws.Range("N136").Value = -14528.5716

$ws = $wb.Worksheets.Item("CUL")
# Row 26: A Grape Idea
$ws.Range("H26").Value = 22223152
$ws.Range("I26").Value = 621.6667
$ws.Range("J26").Value = 33334416
$ws.Range("K26").Value = 1865.0001
$ws.Range("L26").Value = 100003248
$ws.Range("M26").Value = -1577.0001
$ws.Range("N26").Value = -100003824

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 37939370
$ws.Range("I137").Value = 83336140
$ws.Range("J137").Value = 27851194
$ws.Range("K137").Value = 250008420
$ws.Range("L137").Value = 83553582
$ws.Range("M137").Value = -250003320
$ws.Range("N137").Value = -83563782

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 498.36365
$ws.Range("I107").Value = 498.36365
$ws.Range("K107").Value = 498.36365
$ws.Range("M107").Value = 1421.63635

$ws = $wb.Worksheets.Item("LTW")
# Row 18: Simply the Best (M18 previously absent -> now populated)
$ws.Range("H18").Value = 9333.333000000001
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 13000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = -1828
$ws.Range("N18").Value = -13344

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 1537.5
$ws.Range("I82").Value = 2160
$ws.Range("J82").Value = 500
$ws.Range("K82").Value = 2160
$ws.Range("L82").Value = 500
$ws.Range("M82").Value = -1799
$ws.Range("N82").Value = -1222

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 1537.5
$ws.Range("I85").Value = 2160
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 2160
$ws.Range("L85").Value = 500
$ws.Range("M85").Value = -912
$ws.Range("N85").Value = -2996

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1811.4546
$ws.Range("I93").Value = 1746.5714
$ws.Range("J93").Value = 1925
$ws.Range("K93").Value = 1746.5714
$ws.Range("L93").Value = 1925
$ws.Range("M93").Value = -498.5714
$ws.Range("N93").Value = -4421

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 1977.7826
$ws.Range("I132").Value = 1904.4242
$ws.Range("J132").Value = 2164
$ws.Range("K132").Value = 5713.2726
$ws.Range("L132").Value = 6492
$ws.Range("M132").Value = -3183.2726
$ws.Range("N132").Value = -11552

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 9157
$ws.Range("I136").Value = 13274.75
$ws.Range("K136").Value = 39824.25
$ws.Range("M136").Value = -37274.25

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke (L62/M62 removed)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0

# Row 65: Desperate for Diversionaries (L) (L65/M65 removed)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1946.1538
$ws.Range("I81").Value = 1491.6666
$ws.Range("J81").Value = 2335.7144
$ws.Range("K81").Value = 2983.3332
$ws.Range("L81").Value = 4671.4288
$ws.Range("M81").Value = -1922.3332
$ws.Range("N81").Value = -6793.4288

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1946.1538
$ws.Range("I84").Value = 1491.6666
$ws.Range("J84").Value = 2335.7144
$ws.Range("K84").Value = 14916.666
$ws.Range("L84").Value = 23357.144
$ws.Range("M84").Value = -9612.666000000001
$ws.Range("N84").Value = -33965.144

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1853.7715
$ws.Range("I132").Value = 1666.2
$ws.Range("J132").Value = 2322.7
$ws.Range("K132").Value = 4998.6
$ws.Range("L132").Value = 6968.099999999999
$ws.Range("M132").Value = -2468.6
$ws.Range("N132").Value = -12028.1

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 6604.16
$ws.Range("I136").Value = 7019.2383
$ws.Range("J136").Value = 4425
$ws.Range("K136").Value = 21057.7149
$ws.Range("L136").Value = 13275
$ws.Range("M136").Value = -18507.7149
$ws.Range("N136").Value = -18375
